# Append a new data row (row 99) to Sheet 1, mirroring the existing rows
# written by the upstream R script (date, volume, high, low, open, close,
# adj_close-as-text, ticker).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$row = 99
$prev = $row - 1

# Mirror the date cell's style (numFmt "yyyy-mm-dd hh:mm:ss") from the row above
# instead of assigning a NumberFormat string directly (which would mint a
# brand-new style index rather than reusing the existing one).
$ws.Cells.Item($prev, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item($row, 1).Value = 45464.2916666667
$ws.Cells.Item($row, 2).Value = 4720
$ws.Cells.Item($row, 3).Value = 0.689999997615814
$ws.Cells.Item($row, 4).Value = 0.670000016689301
$ws.Cells.Item($row, 5).Value = 0.689999997615814
$ws.Cells.Item($row, 6).Value = 0.670000016689301

# adj_close (G) is written by the source R script as literal text (it matches
# the "close" column's digits exactly but is a string, not a number) - force
# text entry via a temporary "@" number format so the numeric-looking value
# isn't auto-coerced to a number, then drop the format back off the cell so
# it ends up with no explicit style, matching the other data rows.
$gCell = $ws.Cells.Item($row, 7)
$gCell.NumberFormat = "@"
$gCell.Value = "0.670000016689301"
$gCell.Style = "Normal"

$ws.Cells.Item($row, 8).Value = "BWZ.MI"
